# edit.ps1
# Implements the commit: "feat: add 2022-Q1 data"
#
# 1. Adds a new worksheet "2022-Q1" (positioned immediately before "总计")
#    populated with 27 fund rows (same column layout as the existing
#    "2021-Q4" sheet: A index, B code, C name, D size, E position, F ratio,
#    G value, H rank).
# 2. Inserts a new leading row into the "总计" (Total) summary sheet for the
#    "2022-Q1" quarter (27 funds, 24.25 亿元) and renumbers the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by duplicating "2021-Q4" (same header
# labels/styles/column layout) immediately before the "总计" sheet, then
# overwrite its data with the 2022-Q1 numbers.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($totalSheet)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Tab-separated fund rows: A, B(code), C(name), D(scale), E(position), F(ratio), G(value), H(rank)
$fundData = @"
0	000055	广发纳斯达克100指数(QDII) A 美元现汇	75.36	85.84	3.51	2.6451	7
1	270042	广发纳斯达克100指数QDII A	75.36	85.84	3.51	2.6451	7
2	006479	广发纳斯达克100指数（QDII）C人民币	75.36	85.84	3.51	2.6451	7
3	006480	广发纳斯达克100指数（QDII）C美元现汇	75.36	85.84	3.51	2.6451	7
4	513100	国泰纳斯达克100 (QDII-ETF)	41.86	90.74	3.33	1.3939	7
5	513500	博时标普500ETF(QDII)	70.03	90.45	1.98	1.3866	5
6	000988	嘉实全球互联网股票 - 人民币QDII	13.21	85.88	9.69	1.2800	6
7	000989	嘉实全球互联网股票 - 美元现汇QDII	13.21	85.88	9.69	1.2800	6
8	000990	嘉实全球互联网股票 - 美元现钞QDII	13.21	85.88	9.69	1.2800	6
9	000043	嘉实美国成长股票(QDII) -人民币	14.64	94.24	6.12	0.8960	3
10	000044	嘉实美国成长股票(QDII) - 美元现汇	14.64	94.24	6.12	0.8960	3
11	040046	华安纳斯达克100指数QDII - 人民币	22.85	90.93	3.35	0.7655	7
12	040047	华安纳斯达克100指数QDII - 美元现钞	22.85	90.93	3.35	0.7655	7
13	040048	华安纳斯达克100指数QDII - 美元现汇	22.85	90.93	3.35	0.7655	7
14	160213	国泰纳斯达克100指数(QDII)	15.88	90.49	3.34	0.5304	7
15	000834	大成纳斯达克100指数 (QDII)	14.15	89.27	3.28	0.4641	7
16	159941	广发纳斯达克100ETFQDII	11.87	90.26	3.41	0.4048	7
17	003722	易方达纳斯达克100指数美元（QDII-LOF）	9.07	91.29	3.36	0.3048	7
18	161130	易方达纳斯达克100指数人民币（QDII-LOF）	9.07	91.29	3.36	0.3048	7
19	513300	华夏纳斯达克100ETF（QDII）	6.66	93.54	3.44	0.2291	2
20	100055	富国全球科技互联网股票(QDII)	3.01	70.87	4.30	0.1294	4
21	006792	鹏华香港美国互联网股票（LOF）美元现汇	1.43	83.72	7.41	0.1060	2
22	160644	鹏华香港美国互联网股票（LOF）人民币	1.43	83.72	7.41	0.1060	2
23	003718	易方达标普500指数(QDII-LOF) 美元	5.22	91.11	1.98	0.1034	5
24	161125	易方达标普500指数(QDII-LOF) 人民币	5.22	91.11	1.98	0.1034	5
25	012924	华夏新时代灵活配置混合（QDII）美元现汇	2.56	84.71	3.39	0.0868	5
26	012925	华夏新时代灵活配置混合（QDII）美元现钞	2.56	84.71	3.39	0.0868	5
"@

# Make sure there are enough styled rows (A:H) to hold all 27 records; the
# template only had 24 data rows, so clone the formatting of the last
# existing row down across the extra rows before filling in values.
$q1.Range("A24:H24").Copy()
$q1.Range("A25:H28").PasteSpecial(-4122)

$lines = $fundData -split "`n"
$r = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $fields = $line -split "`t"

    $q1.Cells.Item($r, 1).Value = [int]$fields[0]

    # B..G are stored as literal text in the source data (fund codes keep
    # leading zeros, decimal values keep trailing zeros) - force text via a
    # temporary "@" number format, then clear the format again so the saved
    # cell carries no explicit style, matching the other data sheets.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q1.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $fields[$c - 1]
        $cell.ClearFormats()
    }

    $q1.Cells.Item($r, 8).Value = [int]$fields[7]

    $r++
}

# ---------------------------------------------------------------------------
# Step 2: add the new "2022-Q1" line to the "总计" summary sheet and
# renumber the index column (A) sequentially from 0.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Existing rows 2-4 (2021-Q4, 2021-Q2, 2021-Q1) shift down to rows 3-5; clone
# the index-column formatting from row 4 down into the newly used row 5.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = $total.Cells.Item(4, 2).Value
$total.Cells.Item(5, 3).Value = $total.Cells.Item(4, 3).Value
$total.Cells.Item(5, 4).Value = $total.Cells.Item(4, 4).Value

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = $total.Cells.Item(3, 2).Value
$total.Cells.Item(4, 3).Value = $total.Cells.Item(3, 3).Value
$total.Cells.Item(4, 4).Value = $total.Cells.Item(3, 4).Value

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = $total.Cells.Item(2, 2).Value
$total.Cells.Item(3, 3).Value = $total.Cells.Item(2, 3).Value
$total.Cells.Item(3, 4).Value = $total.Cells.Item(2, 4).Value

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 27
$total.Cells.Item(2, 4).Value = 24.25

Write-Output "2022-Q1 sheet added; 总计 sheet updated"
